# Queue , Stack , Extent Reports , CrossBrowser
#
# Adds two new worksheets (QueuePage, Stackpage) to the dsAlgo workbook,
# positioned:
#   ... DataStructure, QueuePage, LinkedlistPage, TreePage, Stackpage
# and tweaks the DataStructure sheet's selection/active-tab state.

$wb = $excel.ActiveWorkbook

# Reusable snippets of "python code" / "run result" shared strings already
# present in the workbook (used verbatim so the engine de-dupes them against
# the existing shared-string table entries instead of creating new ones).
$codeOk = "num1 = 1.5`nnum2 = 6.3`nsum = num1 + num2 `nprint('The sum of {0} and {1} is {2}'.format(num1, num2, sum))`n"
$codeErr = "num1 = 1.5`nnum2 = 6.3`nsum = num1 + num2 `nprint('The sum of {0} and {1} is {2}'.format(num1, num2, `n"
$resultOk = "The sum of 1.5 and 6.3 is 7.8"
$resultErr = "SyntaxError: EOF in multi-line statement on line 6"

# ---------------------------------------------------------------------
# 1. Create "QueuePage" right after "DataStructure"
# ---------------------------------------------------------------------
$dataStructureSheet = $wb.Worksheets.Item("DataStructure")
$queuePage = $wb.Worksheets.Add($null, $dataStructureSheet)
$queuePage.Name = "QueuePage"

$queueLinks = @(
    "ImplementationofQueueinPython",
    "Implementationusingcollectionsdeque",
    "Implementationusingarray ",
    "QueueOperations "
)

$queuePage.Cells.Item(1, 1).Value = "Link"
$queuePage.Cells.Item(1, 2).Value = "Pythoncode"
$queuePage.Cells.Item(1, 3).Value = "Runresult"

$row = 2
foreach ($link in $queueLinks) {
    $queuePage.Cells.Item($row, 1).Value = $link
    $queuePage.Cells.Item($row, 2).Value = $codeOk
    $queuePage.Cells.Item($row, 2).WrapText = $true
    $queuePage.Cells.Item($row, 3).Value = $resultOk
    $row++

    $queuePage.Cells.Item($row, 1).Value = $link
    $queuePage.Cells.Item($row, 2).Value = $codeErr
    $queuePage.Cells.Item($row, 2).WrapText = $true
    $queuePage.Cells.Item($row, 3).Value = $resultErr
    $row++
}

for ($r = 1; $r -le 9; $r++) {
    $queuePage.Rows.Item($r).RowHeight = 22
}

$queuePage.Columns.Item(1).ColumnWidth = 32.166666666666664
$queuePage.Columns.Item(2).ColumnWidth = 10.998697916666666
$queuePage.Columns.Item(3).ColumnWidth = 41.998697916666664

# ---------------------------------------------------------------------
# 2. Create "Stackpage" right after "TreePage"
# ---------------------------------------------------------------------
$treePageSheet = $wb.Worksheets.Item("TreePage")
$stackPage = $wb.Worksheets.Add($null, $treePageSheet)
$stackPage.Name = "Stackpage"

$stackLinks = @(
    "Operationsinstack",
    "Implementation",
    "Applications "
)

$stackPage.Cells.Item(1, 1).Value = "Link"
$stackPage.Cells.Item(1, 2).Value = "Pythoncode"
$stackPage.Cells.Item(1, 3).Value = "Runresult"

$row = 2
foreach ($link in $stackLinks) {
    $stackPage.Cells.Item($row, 1).Value = $link
    $stackPage.Cells.Item($row, 2).Value = $codeOk
    $stackPage.Cells.Item($row, 2).WrapText = $true
    $stackPage.Cells.Item($row, 3).Value = $resultOk
    $row++

    $stackPage.Cells.Item($row, 1).Value = $link
    $stackPage.Cells.Item($row, 2).Value = $codeErr
    $stackPage.Cells.Item($row, 2).WrapText = $true
    $stackPage.Cells.Item($row, 3).Value = $resultErr
    $row++
}

for ($r = 1; $r -le 7; $r++) {
    $stackPage.Rows.Item($r).RowHeight = 61
}

$stackPage.Columns.Item(1).ColumnWidth = 14.998697916666666
$stackPage.Columns.Item(3).ColumnWidth = 41.998697916666664

# ---------------------------------------------------------------------
# 3. View/selection tweaks
# ---------------------------------------------------------------------

# DataStructure: selection becomes B2:C3 (no longer the active tab)
$dataStructureSheet.Activate() | Out-Null
$dataStructureSheet.Range("B2:C3").Select() | Out-Null

# QueuePage: selection becomes B1:C6
$queuePage.Activate() | Out-Null
$queuePage.Range("B1:C6").Select() | Out-Null

# Stackpage: selection becomes A8:XFD24 (whole rows 8-24), and it ends up
# being the active / tabSelected sheet, scrolled so row 3 is at the top.
$stackPage.Activate() | Out-Null
$stackPage.Rows("8:24").Select() | Out-Null
